# Applies the row-data permutation described in the commit
# "Fruta / hortaliza, semanal" to the active worksheet.
#
# The underlying data rows (2-11) keep their row position, but the
# Fecha/Volumen/Precio.../Origen values are reshuffled between rows,
# as if the weekly price records had been re-sorted/re-matched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to write into each row (2..11), taken from the target
# state of the workbook. Columns not listed here (A,B,C,E-L,Q,T) are
# unchanged.
$rowData = @{
    2  = @{ D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; S = 2400; R = "Región de La Araucanía" }
    3  = @{ D = 44574; M = 200; N = 3000; O = 3000; P = 3000; S = 3000; R = "Región de La Araucanía" }
    4  = @{ D = 44616; M = 200; N = 3200; O = 3200; P = 3200; S = 3200; R = "Región de La Araucanía" }
    5  = @{ D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; S = 3000; R = "Región de O'Higgins" }
    6  = @{ D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; S = 1800; R = "Región de La Araucanía" }
    7  = @{ D = 44175; M = 40;  N = 5000; O = 5000; P = 5000; S = 5000; R = "Provincia de Curicó" }
    8  = @{ D = 44551; M = 120; N = 4500; O = 4500; P = 4500; S = 4500; R = "Región de O'Higgins" }
    9  = @{ D = 44592; M = 5;   N = 7500; O = 7500; P = 7500; S = 7500; R = "Región de La Araucanía" }
    10 = @{ D = 44323; M = 20;  N = 3200; O = 3200; P = 3200; S = 3200; R = "Región de La Araucanía" }
    11 = @{ D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; S = 2800; R = "Región de La Araucanía" }
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]

    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
